$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.595.62'
$ws.Range("E2").Value = '  +0.54%  '

$ws.Range("D3").Value = '1.740.70'
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.9989'
$ws.Range("D4").Style = 'Normal'

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '247.11'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +1.35%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = 'Normal'

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.4925'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +2.68%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.2683'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +0.57%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.06290'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +1.01%  '

$ws.Range("D10").Value = '1.736.17'
$ws.Range("E10").Value = '  +0.34%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.07056'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.86%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '15.77'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.6160'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -0.35%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '4.593'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +1.08%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '77.87'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +1.16%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '0.9990'
$ws.Range("D16").Style = 'Normal'

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.000007383'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +6.57%  '

$ws.Range("D18").Value = '26.598.34'
$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.9988'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '11.56'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -1.33%  '

$ws.Range("D21").Value = '1.959.23'
$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '4.594'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +0.72%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '8.742'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -1.69%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '5.261'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -1.09%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '140.28'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +2.78%  '

$ws.Range("E26").Value = '  +0.85%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '1.419'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +0.98%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '108.52'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +1.88%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '1.770'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.14%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '4.054'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +1.78%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.08084'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +1.09%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '3.733'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +0.06%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '0.04628'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +1.71%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '2.611'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -0.19%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.021'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +3.21%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.6385'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -0.36%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.9021'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -3.92%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '2.035'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +2.72%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '2.404'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -0.18%  '

$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.01505'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +0.17%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '101.85'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -5.03%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '5.426'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -4.33%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.3936'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +0.78%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '6.929'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -0.15%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.1190'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -0.20%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.05400'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +1.27%  '

$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '30.58'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -0.70%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '7.799'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -1.10%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '1.273'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -0.12%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '51.81'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +0.95%  '
